# Update crypto price/volume table with latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.445.00"
$ws.Range("E2").Value = "  +5.04%  "
$ws.Range("D3").Value = "2.054.13"
$ws.Range("E3").Value = "  +3.42%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.651"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.85"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +11.94%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.383"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.48"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0769"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.61%  "
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.924"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.83"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").Value = "2.356.71"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.38"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +24.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.57"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.87%  "
$ws.Range("D18").Value = "2.052.91"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("D19").Value = "37.362.35"
$ws.Range("E19").Value = "  +4.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.66"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.04%  "
$ws.Range("E21").Value = "  +3.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.52"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.78"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.69"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.93"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +8.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.88"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.98"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.125"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +29.28%  "
$ws.Range("E31").Value = "  +7.89%  "
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("E33").Value = "  +9.47%  "
$ws.Range("E34").Value = "  +7.96%  "
$ws.Range("E35").Value = "  +5.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.48"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  +3.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.05"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +14.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.05"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +34.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.103"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +15.36%  "
$ws.Range("E42").Value = "  +4.39%  "
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.63"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +8.13%  "
$ws.Range("E45").Value = "  +5.85%  "
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.35"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.20%  "
$ws.Range("E48").Value = "  +2.82%  "
$ws.Range("D49").Value = "1.416.46"
$ws.Range("E49").Value = "  +3.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.89"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +12.28%  "
$ws.Range("E51").Value = "  +1.80%  "
